$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows("7:7").Insert()
$ws.Rows(7).RowHeight = 21.95
$ws.Range("A7").Value = "Help"
$ws.Range("A7").Font.Bold = $true
$ws.Range("B7").Value = "https://ccdb.esss.lu.se/resources/help/ccdb_conventions.pdf"
$ws.Hyperlinks.Add($ws.Range("B7"), "https://ccdb.esss.lu.se/resources/help/ccdb_conventions.pdf")
